# QCMI-BSA-offset.xlsx: fill in the offset totals (column A/C/E/G/I/K/M on
# row 2), box-border + center/wrap those "total" cells, drop the header
# row's leftover style, widen column E to fit, and size rows 1-2 for the
# thicker boxed look. Also nudges the saved window position / selection
# to match the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- header row (A1:N1): strip the stray cell style so it's unstyled ---
$ws.Range("A1:N1").Style = "Normal"

# --- fill in the previously-zeroed "total" cells in row 2 ---
$ws.Range("A2").Value = 4961272.0199999996
$ws.Range("C2").Value = 14866726.869999999
$ws.Range("E2").Value = 24774372.23
$ws.Range("G2").Value = 34680641.270000003
$ws.Range("I2").Value = 44587621.439999998
$ws.Range("K2").Value = 54494404.460000001
$ws.Range("M2").Value = 64401754.57

# --- box border + vertical-center + wrap text on those same total cells ---
# (done as separate passes over the cell set -- interleaving the three
# property writes per-cell leaves more orphaned intermediate cell formats
# behind in the style table)
$totalCols = @(1, 3, 5, 7, 9, 11, 13)
foreach ($col in $totalCols) {
    $ws.Cells.Item(2, $col).Borders.Weight = -4138
}
foreach ($col in $totalCols) {
    $ws.Cells.Item(2, $col).VerticalAlignment = -4108
}
foreach ($col in $totalCols) {
    $ws.Cells.Item(2, $col).WrapText = $true
}

# --- row heights for the now-bolder-looking header/total rows ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75

# --- column E widened to fit its new long number ---
$ws.Columns.Item(5).ColumnWidth = 10.71

# --- selection + saved window position, matching the author's session ---
$null = $ws.Range("M2").Select()
$excel.ActiveWindow.Left = 13380
$excel.ActiveWindow.Top = 4005
